# Generate Report for Handback
# Refreshes the timestamp / handoff-handback tracking columns for the
# "fd31c7fa-9349-463d-91cb-649c56cef66f.md" row after a new handback cycle.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: bump the "Latest HO Xliff Generate Date" for the
# fd31c7fa... file (row 3) to reflect the newly generated handback report.
$overview.Range("G3").Value = "2016-08-29 16:53:21"

# zh-cn sheet, row 3 (fd31c7fa...): refresh handoff/handback datetimes and
# flag that the file now carries metadata.
$zhcn.Range("H3").Value = "2016-08-29 16:53:16"
$zhcn.Range("K3").Value = "2016-08-29 16:53:32"
$zhcn.Range("O3").Value = "'True"

# de-de sheet, row 3 (fd31c7fa...): refresh handoff/handback datetimes.
$dede.Range("H3").Value = "2016-08-29 16:53:21"
$dede.Range("K3").Value = "2016-08-29 16:53:40"
